$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.881.76'
$ws.Range("E2").Value = '  -0.27%  '

# Row 3
$ws.Range("D3").Value = '1.630.58'
$ws.Range("E3").Value = '  -0.67%  '

# Row 4
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.36'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("D6").ClearFormats()

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.24%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.36'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.50%  '

# Row 9
$ws.Range("E9").Value = '  -0.49%  '

# Row 10
$ws.Range("E10").Value = '  -0.38%  '

# Row 11
$ws.Range("E11").Value = '  -0.31%  '

# Row 12
$ws.Range("D12").Value = '1.859.36'
$ws.Range("E12").Value = '  -0.79%  '

# Row 13
$ws.Range("D13").Value = '1.625.77'
$ws.Range("E13").Value = '  -0.91%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.46%  '

# Row 15
$ws.Range("E15").Value = '  -1.42%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.28'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.34%  '

# Row 17
$ws.Range("D17").Value = '27.862.63'
$ws.Range("E17").Value = '  -0.31%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.62'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.50%  '

# Row 19
$ws.Range("E19").Value = '  +0.88%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0721'
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("E21").Value = '  -0.28%  '

# Row 22
$ws.Range("E22").Value = '  -0.95%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.10'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.60%  '

# Row 24
$ws.Range("E24").Value = '  -2.16%  '

# Row 25
$ws.Range("E25").Value = '  +0.60%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.07%  '

# Row 27
$ws.Range("E27").Value = '  -0.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.54'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.81%  '

# Row 29
$ws.Range("E29").Value = '  -0.20%  '

# Row 30
$ws.Range("E30").Value = '  -0.94%  '

# Row 31
$ws.Range("E31").Value = '  -0.68%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.13%  '

# Row 33
$ws.Range("E33").Value = '  -0.09%  '

# Row 34
$ws.Range("D34").Value = '1.394.06'
$ws.Range("E34").Value = '  -0.93%  '

# Row 35
$ws.Range("E35").Value = '  +0.53%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.27%  '

# Row 37
$ws.Range("E37").Value = '  -1.03%  '

# Row 38
$ws.Range("E38").Value = '  +0.44%  '

# Row 39
$ws.Range("E39").Value = '  -1.19%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.853'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.16%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.01'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.27%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.22%  '

# Row 43
$ws.Range("E43").Value = '  +0.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.81'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.13%  '

# Row 45
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("D46").Value = '1.766.71'
$ws.Range("E46").Value = '  -0.91%  '

# Row 47
$ws.Range("E47").Value = '  -2.64%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.15'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.29%  '

# Row 49
$ws.Range("E49").Value = '  +1.53%  '

# Row 50
$ws.Range("E50").Value = '  -0.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.18%  '
